$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 648.1667
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 831.2222
$ws.Range("K17").Value = 297
$ws.Range("L17").Value = 2493.6666
$ws.Range("M17").Value = -129
$ws.Range("N17").Value = -2829.6666
$ws.Range("H21").Value = 39997.5
$ws.Range("I21").Value = 39997.5
$ws.Range("K21").Value = 39997.5
$ws.Range("M21").Value = -39529.5
$ws.Range("H23").Value = 39997.5
$ws.Range("I23").Value = 39997.5
$ws.Range("K23").Value = 39997.5
$ws.Range("M23").Value = -39763.5
$ws.Range("H48").Value = 150
$ws.Range("I48").Value = 150
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 450
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -158
$ws.Range("N48").Value = ""
$ws.Range("H53").Value = 3094.9412
$ws.Range("I53").Value = 4712.909
$ws.Range("J53").Value = 128.66667
$ws.Range("K53").Value = 4712.909
$ws.Range("L53").Value = 128.66667
$ws.Range("M53").Value = -4075.909
$ws.Range("N53").Value = -1402.66667
$ws.Range("H56").Value = 150
$ws.Range("I56").Value = 150
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 450
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 84
$ws.Range("N56").Value = ""
$ws.Range("H74").Value = 11000.333
$ws.Range("I74").Value = 8666.666999999999
$ws.Range("J74").Value = 13334
$ws.Range("K74").Value = 8666.666999999999
$ws.Range("L74").Value = 13334
$ws.Range("M74").Value = -7730.666999999999
$ws.Range("N74").Value = -15206
$ws.Range("H77").Value = 11000.333
$ws.Range("I77").Value = 8666.666999999999
$ws.Range("J77").Value = 13334
$ws.Range("K77").Value = 43333.335
$ws.Range("L77").Value = 66670
$ws.Range("M77").Value = -38653.335
$ws.Range("N77").Value = -76030
$ws.Range("H80").Value = 413.66666
$ws.Range("J80").Value = 464.63635
$ws.Range("L80").Value = 1393.90905
$ws.Range("N80").Value = -3389.90905
$ws.Range("H83").Value = 413.66666
$ws.Range("J83").Value = 464.63635
$ws.Range("L83").Value = 4181.72715
$ws.Range("N83").Value = -14165.72715
$ws.Range("H116").Value = 4371.364
$ws.Range("I116").Value = 4147.5
$ws.Range("K116").Value = 4147.5
$ws.Range("M116").Value = -705.5
$ws.Range("H137").Value = 100004790
$ws.Range("I137").Value = 125005740
$ws.Range("K137").Value = 375017220
$ws.Range("M137").Value = -375014670
$ws.Range("H138").Value = 2249.6428
$ws.Range("I138").Value = 1862.7778
$ws.Range("J138").Value = 2355.1516
$ws.Range("K138").Value = 5588.3334
$ws.Range("L138").Value = 7065.4548
$ws.Range("M138").Value = -448.3334000000004
$ws.Range("N138").Value = -17345.4548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2395.3845
$ws.Range("I86").Value = 2462.8
$ws.Range("K86").Value = 2462.8
$ws.Range("M86").Value = -1339.8
$ws.Range("H89").Value = 2395.3845
$ws.Range("I89").Value = 2462.8
$ws.Range("K89").Value = 12314
$ws.Range("M89").Value = -6698
$ws.Range("H107").Value = 4899.273
$ws.Range("I107").Value = 4970.778
$ws.Range("J107").Value = 4577.5
$ws.Range("K107").Value = 4970.778
$ws.Range("L107").Value = 4577.5
$ws.Range("M107").Value = -3050.778
$ws.Range("N107").Value = -8417.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4638.0835
$ws.Range("I31").Value = 4952.2
$ws.Range("J31").Value = 4413.7144
$ws.Range("K31").Value = 4952.2
$ws.Range("L31").Value = 4413.7144
$ws.Range("M31").Value = -4657.2
$ws.Range("N31").Value = -5003.7144
$ws.Range("H34").Value = 4638.0835
$ws.Range("I34").Value = 4952.2
$ws.Range("J34").Value = 4413.7144
$ws.Range("K34").Value = 4952.2
$ws.Range("L34").Value = 4413.7144
$ws.Range("M34").Value = -4750.2
$ws.Range("N34").Value = -4817.7144
$ws.Range("H86").Value = 50009316
$ws.Range("I86").Value = 76930800
$ws.Range("K86").Value = 76930800
$ws.Range("M86").Value = -76929677
$ws.Range("H89").Value = 50009316
$ws.Range("I89").Value = 76930800
$ws.Range("K89").Value = 384654000
$ws.Range("M89").Value = -384648384
$ws.Range("H94").Value = 970.4375
$ws.Range("J94").Value = 680.5714
$ws.Range("L94").Value = 680.5714
$ws.Range("N94").Value = -1582.5714
$ws.Range("H99").Value = 14817665
$ws.Range("I99").Value = 6669571
$ws.Range("J99").Value = 55558136
$ws.Range("K99").Value = 6669571
$ws.Range("L99").Value = 55558136
$ws.Range("M99").Value = -6668073
$ws.Range("N99").Value = -55561132
$ws.Range("H126").Value = 14817665
$ws.Range("I126").Value = 6669571
$ws.Range("J126").Value = 55558136
$ws.Range("K126").Value = 20008713
$ws.Range("L126").Value = 166674408
$ws.Range("M126").Value = -20006243
$ws.Range("N126").Value = -166679348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29444620
$ws.Range("J2").Value = 47619304
$ws.Range("L2").Value = 285715824
$ws.Range("N2").Value = -285716050
$ws.Range("H7").Value = 295.68182
$ws.Range("I7").Value = 267.14285
$ws.Range("K7").Value = 801.4285500000001
$ws.Range("M7").Value = -689.4285500000001
$ws.Range("H23").Value = 83
$ws.Range("J23").Value = 83
$ws.Range("L23").Value = 249
$ws.Range("N23").Value = -719
$ws.Range("H39").Value = 2824.75
$ws.Range("J39").Value = 3533
$ws.Range("L39").Value = 10599
$ws.Range("N39").Value = -11187
$ws.Range("H50").Value = 20000476
$ws.Range("J50").Value = 25000218
$ws.Range("L50").Value = 75000654
$ws.Range("N50").Value = -75001616
$ws.Range("H53").Value = 20000476
$ws.Range("J53").Value = 25000218
$ws.Range("L53").Value = 75000654
$ws.Range("N53").Value = -75001616
$ws.Range("H107").Value = 555.7037
$ws.Range("J107").Value = 587.9091
$ws.Range("L107").Value = 1763.7273
$ws.Range("N107").Value = -5603.7273
$ws.Range("H129").Value = 30955054
$ws.Range("J129").Value = 12503873
$ws.Range("L129").Value = 37511619
$ws.Range("N129").Value = -37521619

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6224.75
$ws.Range("J126").Value = 6224.75
$ws.Range("L126").Value = 18674.25
$ws.Range("N126").Value = -23614.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2438.875
$ws.Range("I61").Value = 2802.0908
$ws.Range("K61").Value = 2802.0908
$ws.Range("M61").Value = -2600.0908
$ws.Range("H113").Value = 2438.875
$ws.Range("I113").Value = 2802.0908
$ws.Range("K113").Value = 2802.0908
$ws.Range("M113").Value = -632.0907999999999
$ws.Range("H132").Value = 3323.6667
$ws.Range("I132").Value = 4367
$ws.Range("J132").Value = 2489
$ws.Range("K132").Value = 13101
$ws.Range("L132").Value = 7467
$ws.Range("M132").Value = -10571
$ws.Range("N132").Value = -12527
$ws.Range("H140").Value = 64412.723
$ws.Range("I140").Value = 40000
$ws.Range("J140").Value = 79948.09
$ws.Range("K140").Value = 40000
$ws.Range("L140").Value = 79948.09
$ws.Range("M140").Value = -34820
$ws.Range("N140").Value = -90308.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 669.8
$ws.Range("I29").Value = 606.6667
$ws.Range("J29").Value = 764.5
$ws.Range("K29").Value = 606.6667
$ws.Range("L29").Value = 764.5
$ws.Range("M29").Value = -316.6667
$ws.Range("N29").Value = -1344.5
$ws.Range("H117").Value = 73998
$ws.Range("J117").Value = 73998
$ws.Range("L117").Value = 73998
$ws.Range("N117").Value = -83176
$ws.Range("H136").Value = 60389.47
$ws.Range("I136").Value = 4855.5386
$ws.Range("J136").Value = 240874.75
$ws.Range("K136").Value = 14566.6158
$ws.Range("L136").Value = 722624.25
$ws.Range("M136").Value = -12016.6158
$ws.Range("N136").Value = -727724.25

Write-Host "Applied 204 cell updates across 7 sheets"